$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column J ("Date") for existing rows 2-17, plus new row 18 (new
#    CPU entry) and a trailing, otherwise-empty row 19.
# ---------------------------------------------------------------------

# Header cell J1 - value first, style (blue header look) copied after so
# the shared string for "Date" is minted only once we get there (it must
# be the LAST new shared string, matching the source order).

# Fill the date column for the already-existing rows (2-17) with the
# same serial date value used throughout the original sheet.
$ws.Range("J2:J17").Value = 45323

# New row of data (row 18). Columns are set in a specific order so any
# brand-new shared strings are appended in the same sequence as the
# target workbook (RAM before CPU, then the four timing columns).
$ws.Range("A18").Value = "Windows"
$ws.Range("B18").Value = "PC"
$ws.Range("D18").Value = "3.6"
$ws.Range("E18").Value = "DDR3 16GB"
$ws.Range("C18").Value = "AMD Phenom II X4 975"
$ws.Range("F18").Value = "48ms"
$ws.Range("G18").Value = "477ms"
$ws.Range("H18").Value = "4.8s"
$ws.Range("I18").Value = "48s"

# Header for the new column, after the row-18 strings so "Date" becomes
# the final newly-minted shared string.
$ws.Range("J1").Value = "Date"

# Date value for the new row, and an otherwise-blank row 19 underneath
# it (only the date column is touched there, and it stays valueless -
# formatted but empty, same as the source).
$ws.Range("J18").Value = 45324

# ---------------------------------------------------------------------
# 2. Number formatting - stamp a single cell with a date format that
#    resolves to the built-in numFmtId 14, then fan that exact style out
#    to the rest of the date column via a format-only paste so every
#    cell shares one cellXfs entry instead of minting a new one each.
# ---------------------------------------------------------------------
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").Copy()
$ws.Range("J3:J18").PasteSpecial(-4122)
$ws.Range("J19").PasteSpecial(-4122)

# Give the new header cell the same look (fill/border/font/alignment) as
# the rest of row 1's header band.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Re-apply the header text/value - PasteSpecial(formats) shouldn't touch
# it, but make sure nothing got clobbered.
$ws.Range("J1").Value = "Date"

# Clear the leftover marching-ants selection from the copy operations.
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Column widths - the new Date column, and tiny re-flow nudges on the
#    existing columns (best effort given the host's width quantisation).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 23.529947916666668
$ws.Columns.Item(3).ColumnWidth = 19.799479166666668
$ws.Columns.Item(4).ColumnWidth = 4.983072916666667
$ws.Columns.Item(6).ColumnWidth = 9.436197916666666
$ws.Columns.Item(7).ColumnWidth = 10.709635416666666
$ws.Columns.Item(8).ColumnWidth = 11.619791666666666
$ws.Columns.Item(9).ColumnWidth = 12.436197916666666
$ws.Columns.Item(10).ColumnWidth = 13.436197916666666

# ---------------------------------------------------------------------
# 4. View state - zoom level and active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("L23").Select()
